{"js": "// The document body contains a single 20x5 table of simple arithmetic\n// expressions (\"85-19=\", etc.). The edit replaces the text of every one\n// of the 100 cells with a new expression, in row-major (reading) order,\n// while leaving all paragraph/run formatting untouched.\n//\n// `pairs` holds [oldText, newText] for each cell in document order, taken\n// directly from the canonical OOXML diff. We apply replacements\n// positionally (not via a global search/replace) because a few of the\n// original expressions repeat verbatim elsewhere in the table but map to\n// different replacement text depending on where they sit.\nconst pairs = [[\"45-22=\", \"85-19=\"], [\"25+3=\", \"59-12=\"], [\"40-14=\", \"62+31=\"], [\"85-16=\", \"76-67=\"], [\"55+5=\", \"78-23=\"], [\"38+51=\", \"42+23=\"], [\"5-1=\", \"92-41=\"], [\"26-22=\", \"98-55=\"], [\"7+2=\", \"2+96=\"], [\"17+82=\", \"60+10=\"], [\"29-7=\", \"88-85=\"], [\"8+45=\", \"91-62=\"], [\"51-15=\", \"38-8=\"], [\"50-41=\", \"16+81=\"], [\"63+17=\", \"59+27=\"], [\"31+64=\", \"71-29=\"], [\"90-40=\", \"99-59=\"], [\"88-53=\", \"10+41=\"], [\"51-18=\", \"13+5=\"], [\"34-4=\", \"19+4=\"], [\"8+34=\", \"13+53=\"], [\"65+13=\", \"75-56=\"], [\"48-36=\", \"81-54=\"], [\"26+54=\", \"36-18=\"], [\"26-13=\", \"80-46=\"], [\"39-1=\", \"63-12=\"], [\"34+29=\", \"23+24=\"], [\"12-0=\", \"26+51=\"], [\"68+2=\", \"48+2=\"], [\"78-62=\", \"17+14=\"], [\"68-45=\", \"61-16=\"], [\"91-70=\", \"57+12=\"], [\"26+3=\", \"85+4=\"], [\"97-32=\", \"88-57=\"], [\"67-59=\", \"84+11=\"], [\"36+43=\", \"60+27=\"], [\"24+69=\", \"46-4=\"], [\"84+1=\", \"46-43=\"], [\"92-38=\", \"73+23=\"], [\"35-12=\", \"61-55=\"], [\"30+44=\", \"2+55=\"], [\"34+40=\", \"1+54=\"], [\"96-78=\", \"0+23=\"], [\"75+6=\", \"94-27=\"], [\"72-31=\", \"62-32=\"], [\"78+17=\", \"73-30=\"], [\"4+89=\", \"73-40=\"], [\"38+14=\", \"91+1=\"], [\"59+7=\", \"92-3=\"], [\"89-15=\", \"53-34=\"], [\"51+28=\", \"5+35=\"], [\"65-47=\", \"1+28=\"], [\"28+30=\", \"11+35=\"], [\"4+6=\", \"14+23=\"], [\"47+52=\", \"8+4=\"], [\"36-14=\", \"66-44=\"], [\"14+43=\", \"14+50=\"], [\"3+0=\", \"71+11=\"], [\"33+4=\", \"83-74=\"], [\"10+12=\", \"54-7=\"], [\"80+15=\", \"19-5=\"], [\"47+46=\", \"67-52=\"], [\"9+71=\", \"24+8=\"], [\"63-17=\", \"12+22=\"], [\"66-46=\", \"0+9=\"], [\"38+35=\", \"73-26=\"], [\"53+23=\", \"12-2=\"], [\"5+19=\", \"64-33=\"], [\"52+31=\", \"75+21=\"], [\"94-35=\", \"22-19=\"], [\"9+1=\", \"18-4=\"], [\"31+51=\", \"94+3=\"], [\"35+53=\", \"74-49=\"], [\"15+80=\", \"12+80=\"], [\"7+4=\", \"54-51=\"], [\"1+63=\", \"99-23=\"], [\"69-53=\", \"70-1=\"], [\"96-76=\", \"62+32=\"], [\"54-26=\", \"39-14=\"], [\"19+77=\", \"0+38=\"], [\"84-35=\", \"21+65=\"], [\"26-22=\", \"51+10=\"], [\"86+11=\", \"84-15=\"], [\"66+4=\", \"47-19=\"], [\"31+47=\", \"29+27=\"], [\"44+42=\", \"73-12=\"], [\"48+39=\", \"50+45=\"], [\"2+68=\", \"54+16=\"], [\"70-16=\", \"58-29=\"], [\"57-22=\", \"22+1=\"], [\"23+11=\", \"22+17=\"], [\"95-17=\", \"18+65=\"], [\"69-41=\", \"67-18=\"], [\"81+12=\", \"40-12=\"], [\"99-85=\", \"31+20=\"], [\"43+12=\", \"1+47=\"], [\"19+36=\", \"44-18=\"], [\"12+7=\", \"69-14=\"], [\"63-21=\", \"65-49=\"], [\"26+44=\", \"21+77=\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nconst grid = table.values; // string[][], row-major\n\nlet idx = 0;\nfor (let r = 0; r < grid.length; r++) {\n  for (let c = 0; c < grid[r].length; c++) {\n    if (idx >= pairs.length) {\n      break;\n    }\n    const [, newText] = pairs[idx];\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    // Replacing the paragraph's text (rather than clearing the cell body\n    // and inserting a new run) keeps the existing run's rPr/pPr\n    // (font, size, alignment) intact.\n    para.insertText(newText, Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document body contains a single 20x5 table of simple arithmetic\n# expressions (\"85-19=\", etc.). The edit replaces the text of every one\n# of the 100 cells with a new expression, in row-major (reading) order,\n# while leaving all paragraph/run formatting untouched.\n#\n# $pairs holds (oldText, newText) for each cell in document order, taken\n# directly from the canonical OOXML diff. Replacements are applied\n# positionally (row by row, left to right) rather than via a global\n# search/replace, because a few of the original expressions repeat\n# verbatim elsewhere in the table but map to different replacement text\n# depending on where they sit.\n$pairs = @(\n    @('45-22=', '85-19='),\n    @('25+3=', '59-12='),\n    @('40-14=', '62+31='),\n    @('85-16=', '76-67='),\n    @('55+5=', '78-23='),\n    @('38+51=', '42+23='),\n    @('5-1=', '92-41='),\n    @('26-22=', '98-55='),\n    @('7+2=', '2+96='),\n    @('17+82=', '60+10='),\n    @('29-7=', '88-85='),\n    @('8+45=', '91-62='),\n    @('51-15=', '38-8='),\n    @('50-41=', '16+81='),\n    @('63+17=', '59+27='),\n    @('31+64=', '71-29='),\n    @('90-40=', '99-59='),\n    @('88-53=', '10+41='),\n    @('51-18=', '13+5='),\n    @('34-4=', '19+4='),\n    @('8+34=', '13+53='),\n    @('65+13=', '75-56='),\n    @('48-36=', '81-54='),\n    @('26+54=', '36-18='),\n    @('26-13=', '80-46='),\n    @('39-1=', '63-12='),\n    @('34+29=', '23+24='),\n    @('12-0=', '26+51='),\n    @('68+2=', '48+2='),\n    @('78-62=', '17+14='),\n    @('68-45=', '61-16='),\n    @('91-70=', '57+12='),\n    @('26+3=', '85+4='),\n    @('97-32=', '88-57='),\n    @('67-59=', '84+11='),\n    @('36+43=', '60+27='),\n    @('24+69=', '46-4='),\n    @('84+1=', '46-43='),\n    @('92-38=', '73+23='),\n    @('35-12=', '61-55='),\n    @('30+44=', '2+55='),\n    @('34+40=', '1+54='),\n    @('96-78=', '0+23='),\n    @('75+6=', '94-27='),\n    @('72-31=', '62-32='),\n    @('78+17=', '73-30='),\n    @('4+89=', '73-40='),\n    @('38+14=', '91+1='),\n    @('59+7=', '92-3='),\n    @('89-15=', '53-34='),\n    @('51+28=', '5+35='),\n    @('65-47=', '1+28='),\n    @('28+30=', '11+35='),\n    @('4+6=', '14+23='),\n    @('47+52=', '8+4='),\n    @('36-14=', '66-44='),\n    @('14+43=', '14+50='),\n    @('3+0=', '71+11='),\n    @('33+4=', '83-74='),\n    @('10+12=', '54-7='),\n    @('80+15=', '19-5='),\n    @('47+46=', '67-52='),\n    @('9+71=', '24+8='),\n    @('63-17=', '12+22='),\n    @('66-46=', '0+9='),\n    @('38+35=', '73-26='),\n    @('53+23=', '12-2='),\n    @('5+19=', '64-33='),\n    @('52+31=', '75+21='),\n    @('94-35=', '22-19='),\n    @('9+1=', '18-4='),\n    @('31+51=', '94+3='),\n    @('35+53=', '74-49='),\n    @('15+80=', '12+80='),\n    @('7+4=', '54-51='),\n    @('1+63=', '99-23='),\n    @('69-53=', '70-1='),\n    @('96-76=', '62+32='),\n    @('54-26=', '39-14='),\n    @('19+77=', '0+38='),\n    @('84-35=', '21+65='),\n    @('26-22=', '51+10='),\n    @('86+11=', '84-15='),\n    @('66+4=', '47-19='),\n    @('31+47=', '29+27='),\n    @('44+42=', '73-12='),\n    @('48+39=', '50+45='),\n    @('2+68=', '54+16='),\n    @('70-16=', '58-29='),\n    @('57-22=', '22+1='),\n    @('23+11=', '22+17='),\n    @('95-17=', '18+65='),\n    @('69-41=', '67-18='),\n    @('81+12=', '40-12='),\n    @('99-85=', '31+20='),\n    @('43+12=', '1+47='),\n    @('19+36=', '44-18='),\n    @('12+7=', '69-14='),\n    @('63-21=', '65-49='),\n    @('26+44=', '21+77=')\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -ge $pairs.Count) { break }\n        $newText = $pairs[$idx][1]\n\n        $cellRange = $t.Cell($r, $c).Range\n        # Trim the table cell's trailing end-of-cell mark (one \"character\"\n        # move covers the \\r\\a pair) so we only overwrite the visible text,\n        # keeping the run's formatting (font/size) and the cell mark intact.\n        $cellRange.MoveEnd(1, -1) | Out-Null\n        $cellRange.Text = $newText\n\n        $idx++\n    }\n}\n"}
